# Update "想去人数" (F column) values on both the "展览" and "全部类型"
# worksheets. Both sheets contain identical data, so the same set of
# row/value updates is applied to each sheet.

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 2192
    5  = 12987
    8  = 512
    9  = 476
    11 = 971
    12 = 13732
    13 = 14258
    18 = 31
    25 = 5333
    26 = 936
    27 = 12
    28 = 291
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
